# Fixing new UTM origin
# Moves the existing ("initial") area/UTM-boundary table down by two rows,
# labels it, and adds a corrected ("new") table alongside it with updated
# UTM-low/UTM-high/zone values for the first few areas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing table (B1:F6) down two rows -> becomes B3:F8,
# freeing up row 2 for the "initial" label and leaving row 1 empty
# (matches the new dimension B2:L8 / header row now at row 3).
$ws.Rows("1:2").Insert()

# Label the original table as "initial" and give it the accent fill.
$ws.Range("B2").Value = "initial"

# Build the corrected ("new") table at H3:L8 by copying the header +
# data from the original table (now at B3:F8) and then overwriting the
# cells whose values actually changed.
$ws.Range("B3:F8").Copy()
$ws.Range("H3").PasteSpecial()
$ws.Range("H2").Value = "new "

# --- corrected values for the "new" table ---
# area 1: UTM-low 47->48, UTM-high 51->50, zone-low T->U
$ws.Range("I4").Value = 48
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = "U"

# area 2: UTM-low 51->53, zone-high V->U
$ws.Range("I5").Value = 53
$ws.Range("L5").Value = "U"

# area 3: UTM-low 52->53, zone-high X->W
$ws.Range("I6").Value = 53
$ws.Range("L6").Value = "W"

# areas 4 and 5 are unchanged between "initial" and "new".

# --- formatting: accent-colored fills (theme Accent1, with tints) ---
# Row 2 labels ("initial" / "new") -> solid Accent1, no tint
$ws.Range("B2").Interior.ThemeColor = 5
$ws.Range("B2").Interior.TintAndShade = 0
$ws.Range("H2").Interior.ThemeColor = 5
$ws.Range("H2").Interior.TintAndShade = 0

# Header rows (row 3) for both tables -> Accent1, Lighter 40%
$ws.Range("B3:F3").Interior.ThemeColor = 5
$ws.Range("B3:F3").Interior.TintAndShade = 0.39997558519241921
$ws.Range("H3:L3").Interior.ThemeColor = 5
$ws.Range("H3:L3").Interior.TintAndShade = 0.39997558519241921

# Index columns (B and H, rows 4-8) for both tables -> Accent1, Lighter 60%
$ws.Range("B4:B8").Interior.ThemeColor = 5
$ws.Range("B4:B8").Interior.TintAndShade = 0.59999389629810485
$ws.Range("H4:H8").Interior.ThemeColor = 5
$ws.Range("H4:H8").Interior.TintAndShade = 0.59999389629810485

# Keep the active selection near the new table, like the source file.
$ws.Range("L17").Select()
